$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that currently sits inside the
#    "No lines of more than 2 collinear points" test case (between the
#    " (6,7" run and the closing ")" run).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Fix test case 3's point list:
#    "Points: (0,0) (1,1) (2,2) (2,3) (0,2) (10,5) (5,3) (2,4) (2,8) (2,1)"
#    becomes
#    "Points: (0,0) (1,1) (2,2) (2,3) (0,2) (10,5) (5,3) (2,4) (3,8) (5,2)"
#    (point (2,8) -> (3,8) and point (2,1) -> (5,2); the old set had a
#    group of 5 collinear points, see issue #16).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Points: (0,0) (1,1) (2,2) (2,3) (0,2) (10,5) (5,3) (2,4) (2,8) (2,1)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base = $rng.Start

# Replace the middle chunk of the list (toggling formatting forces the
# editor to keep this as its own run instead of re-merging it with the
# untouched text around it).
$part1 = $d.Range($base + 27, $base + 65)
$part1.Bold = $true
$part1.Text = "2,3) (0,2) (10,5) (5,3) (2,4) (3,8) (5"
$part1.Bold = $false

# Replace the final digit of the last point.
$part2 = $d.Range($base + 66, $base + 67)
$part2.Bold = $true
$part2.Text = "2"
$part2.Bold = $false

# ------------------------------------------------------------------
# 3) Re-insert the _GoBack bookmark right after the closing ")" of the
#    point list we just edited (collapsed/zero-length bookmark).
# ------------------------------------------------------------------
$endRng = $d.Range($base + 68, $base + 68)
$d.Bookmarks.Add("_GoBack", $endRng)
